$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 640.8421
$ws.Range("I2").Value = 128.84616
$ws.Range("K2").Value = 128.84616
$ws.Range("M2").Value = -15.84616
$ws.Range("H61").Value = 999.5
$ws.Range("I61").Value = 999.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2998.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2826.5
$ws.Range("N61").Value = $null
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null
$ws.Range("H80").Value = 760.55554
$ws.Range("I80").Value = 350
$ws.Range("J80").Value = 1089
$ws.Range("K80").Value = 1050
$ws.Range("L80").Value = 3267
$ws.Range("M80").Value = -52
$ws.Range("N80").Value = -5263
$ws.Range("H83").Value = 760.55554
$ws.Range("I83").Value = 350
$ws.Range("J83").Value = 1089
$ws.Range("K83").Value = 3150
$ws.Range("L83").Value = 9801
$ws.Range("M83").Value = 1842
$ws.Range("N83").Value = -19785
$ws.Range("H92").Value = 669.0952
$ws.Range("I92").Value = 580.0526
$ws.Range("K92").Value = 580.0526
$ws.Range("M92").Value = 667.9474
$ws.Range("H93").Value = 28749.5
$ws.Range("J93").Value = 28749.5
$ws.Range("L93").Value = 28749.5
$ws.Range("N93").Value = -33741.5
$ws.Range("H97").Value = 1258.4
$ws.Range("J97").Value = 1258.4
$ws.Range("L97").Value = 3775.2
$ws.Range("N97").Value = -4767.200000000001
$ws.Range("H101").Value = 999.2
$ws.Range("I101").Value = 999
$ws.Range("J101").Value = 999.5
$ws.Range("K101").Value = 2997
$ws.Range("L101").Value = 2998.5
$ws.Range("M101").Value = -1375
$ws.Range("N101").Value = -6242.5
$ws.Range("H103").Value = 478
$ws.Range("I103").Value = 217
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 651
$ws.Range("L103").Value = 3000
$ws.Range("M103").Value = -65
$ws.Range("N103").Value = -4172
$ws.Range("H116").Value = 17044.53
$ws.Range("I116").Value = 17134.363
$ws.Range("J116").Value = 16879.834
$ws.Range("K116").Value = 17134.363
$ws.Range("L116").Value = 16879.834
$ws.Range("M116").Value = -13692.363
$ws.Range("N116").Value = -23763.834
$ws.Range("H132").Value = 2871.5356
$ws.Range("I132").Value = 2860.7917
$ws.Range("K132").Value = 8582.375100000001
$ws.Range("M132").Value = -6052.375100000001
$ws.Range("H134").Value = 57023.285
$ws.Range("J134").Value = 57023.285
$ws.Range("L134").Value = 57023.285
$ws.Range("N134").Value = -67163.285
$ws.Range("H137").Value = 1947.1351
$ws.Range("I137").Value = 1850.1613
$ws.Range("K137").Value = 5550.4839
$ws.Range("M137").Value = -3000.4839
$ws.Range("H138").Value = 3975.8733
$ws.Range("I138").Value = 3418.641
$ws.Range("K138").Value = 10255.923
$ws.Range("M138").Value = -5115.923000000001
$ws.Range("H141").Value = 3029.8462
$ws.Range("I141").Value = 1693.8
$ws.Range("K141").Value = 5081.4
$ws.Range("M141").Value = 98.60000000000036

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8320.569
$ws.Range("I32").Value = 3993.2126
$ws.Range("J32").Value = 26810.182
$ws.Range("K32").Value = 3993.2126
$ws.Range("L32").Value = 26810.182
$ws.Range("M32").Value = -3706.2126
$ws.Range("N32").Value = -27384.182
$ws.Range("H45").Value = 3125.6667
$ws.Range("I45").Value = 1334.75
$ws.Range("K45").Value = 1334.75
$ws.Range("M45").Value = -957.75
$ws.Range("H61").Value = 2986.9644
$ws.Range("I61").Value = 1907.8572
$ws.Range("J61").Value = 4066.0715
$ws.Range("K61").Value = 1907.8572
$ws.Range("L61").Value = 4066.0715
$ws.Range("M61").Value = -1695.8572
$ws.Range("N61").Value = -4490.0715
$ws.Range("H96").Value = 108644
$ws.Range("J96").Value = 108644
$ws.Range("L96").Value = 108644
$ws.Range("N96").Value = -114136
$ws.Range("H102").Value = 2439.8
$ws.Range("I102").Value = 2429.3333
$ws.Range("J102").Value = 2455.5
$ws.Range("K102").Value = 2429.3333
$ws.Range("L102").Value = 2455.5
$ws.Range("M102").Value = -807.3332999999998
$ws.Range("N102").Value = -5699.5
$ws.Range("H131").Value = 98000
$ws.Range("J131").Value = 98000
$ws.Range("L131").Value = 98000
$ws.Range("N131").Value = -108080
$ws.Range("H132").Value = 4055.12
$ws.Range("I132").Value = 2835.842
$ws.Range("J132").Value = 7916.1665
$ws.Range("K132").Value = 8507.526
$ws.Range("L132").Value = 23748.4995
$ws.Range("M132").Value = -5977.526
$ws.Range("N132").Value = -28808.4995
$ws.Range("H133").Value = 39252
$ws.Range("J133").Value = 39252
$ws.Range("L133").Value = 39252
$ws.Range("N133").Value = -44312
$ws.Range("H136").Value = 2986.9644
$ws.Range("I136").Value = 1907.8572
$ws.Range("J136").Value = 4066.0715
$ws.Range("K136").Value = 5723.571599999999
$ws.Range("L136").Value = 12198.2145
$ws.Range("M136").Value = -3173.571599999999
$ws.Range("N136").Value = -17298.2145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1277.8214
$ws.Range("I20").Value = 1159.4117
$ws.Range("J20").Value = 1460.8182
$ws.Range("K20").Value = 1159.4117
$ws.Range("L20").Value = 1460.8182
$ws.Range("M20").Value = -912.4117000000001
$ws.Range("N20").Value = -1954.8182
$ws.Range("H21").Value = 34966.75
$ws.Range("J21").Value = 34966.75
$ws.Range("L21").Value = 34966.75
$ws.Range("N21").Value = -35438.75
$ws.Range("H22").Value = 981
$ws.Range("J22").Value = 792
$ws.Range("L22").Value = 792
$ws.Range("N22").Value = -1138
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("H26").Value = 16225
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = $null
$ws.Range("H40").Value = 43985
$ws.Range("J40").Value = 43985
$ws.Range("L40").Value = 43985
$ws.Range("N40").Value = -44515
$ws.Range("H86").Value = 1796.6957
$ws.Range("I86").Value = 1831.6428
$ws.Range("K86").Value = 1831.6428
$ws.Range("M86").Value = -708.6428000000001
$ws.Range("H89").Value = 1796.6957
$ws.Range("I89").Value = 1831.6428
$ws.Range("K89").Value = 9158.214
$ws.Range("M89").Value = -3542.214
$ws.Range("H93").Value = 50444
$ws.Range("J93").Value = 50444
$ws.Range("L93").Value = 50444
$ws.Range("N93").Value = -54188
$ws.Range("H94").Value = 897.9459000000001
$ws.Range("J94").Value = 1411
$ws.Range("L94").Value = 1411
$ws.Range("N94").Value = -2313
$ws.Range("H99").Value = 2125
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996
$ws.Range("H105").Value = 3491.0833
$ws.Range("I105").Value = 3462.15
$ws.Range("K105").Value = 3462.15
$ws.Range("M105").Value = -1715.15
$ws.Range("H107").Value = 2727.2424
$ws.Range("I107").Value = 2768.9355
$ws.Range("K107").Value = 2768.9355
$ws.Range("M107").Value = -848.9355
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null
$ws.Range("H134").Value = 3023.7
$ws.Range("I134").Value = 2137.4443
$ws.Range("K134").Value = 6412.3329
$ws.Range("M134").Value = -3877.3329

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5425.1577
$ws.Range("I31").Value = 2281.8
$ws.Range("K31").Value = 2281.8
$ws.Range("M31").Value = -1986.8
$ws.Range("H34").Value = 5425.1577
$ws.Range("I34").Value = 2281.8
$ws.Range("K34").Value = 2281.8
$ws.Range("M34").Value = -2079.8
$ws.Range("H52").Value = 60738.668
$ws.Range("J52").Value = 58933.5
$ws.Range("L52").Value = 58933.5
$ws.Range("N52").Value = -59521.5
$ws.Range("H58").Value = 4649.65
$ws.Range("I58").Value = 3181.2727
$ws.Range("J58").Value = 6444.3335
$ws.Range("K58").Value = 3181.2727
$ws.Range("L58").Value = 6444.3335
$ws.Range("M58").Value = -2978.2727
$ws.Range("N58").Value = -6850.3335
$ws.Range("H86").Value = 7500
$ws.Range("I86").Value = 7500
$ws.Range("K86").Value = 7500
$ws.Range("M86").Value = -6377
$ws.Range("H89").Value = 7500
$ws.Range("I89").Value = 7500
$ws.Range("K89").Value = 37500
$ws.Range("M89").Value = -31884
$ws.Range("H94").Value = 5182.25
$ws.Range("I94").Value = 5198.5
$ws.Range("J94").Value = 5166
$ws.Range("K94").Value = 5198.5
$ws.Range("L94").Value = 5166
$ws.Range("M94").Value = -4747.5
$ws.Range("N94").Value = -6068
$ws.Range("H105").Value = 4166.6665
$ws.Range("I105").Value = 4166.6665
$ws.Range("K105").Value = 4166.6665
$ws.Range("M105").Value = -2419.6665
$ws.Range("H132").Value = 6080.15
$ws.Range("I132").Value = 4185.143
$ws.Range("K132").Value = 12555.429
$ws.Range("M132").Value = -10025.429
$ws.Range("H134").Value = 18521676
$ws.Range("I134").Value = 27780364
$ws.Range("K134").Value = 83341092
$ws.Range("M134").Value = -83338557
$ws.Range("H136").Value = 4649.65
$ws.Range("I136").Value = 3181.2727
$ws.Range("J136").Value = 6444.3335
$ws.Range("K136").Value = 9543.8181
$ws.Range("L136").Value = 19333.0005
$ws.Range("M136").Value = -6993.8181
$ws.Range("N136").Value = -24433.0005
$ws.Range("H139").Value = 118582.5
$ws.Range("J139").Value = 118582.5
$ws.Range("L139").Value = 118582.5
$ws.Range("N139").Value = -128862.5
$ws.Range("H140").Value = 121298
$ws.Range("J140").Value = 121298
$ws.Range("L140").Value = 121298
$ws.Range("N140").Value = -131658

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = $null
$ws.Range("N36").Value = $null
$ws.Range("H70").Value = 18749.25
$ws.Range("I70").Value = 17748.5
$ws.Range("K70").Value = 53245.5
$ws.Range("M70").Value = -52930.5
$ws.Range("H73").Value = 18749.25
$ws.Range("I73").Value = 17748.5
$ws.Range("K73").Value = 53245.5
$ws.Range("M73").Value = -52153.5
$ws.Range("H75").Value = 3482.2307
$ws.Range("I75").Value = 3233.3333
$ws.Range("J75").Value = 3556.9
$ws.Range("K75").Value = 9699.999899999999
$ws.Range("L75").Value = 10670.7
$ws.Range("M75").Value = -8701.999899999999
$ws.Range("N75").Value = -12666.7
$ws.Range("H76").Value = 14688.909
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 60000
$ws.Range("N76").Value = -60766
$ws.Range("H78").Value = 3482.2307
$ws.Range("I78").Value = 3233.3333
$ws.Range("J78").Value = 3556.9
$ws.Range("K78").Value = 29099.9997
$ws.Range("L78").Value = 32012.1
$ws.Range("M78").Value = -24107.9997
$ws.Range("N78").Value = -41996.10000000001
$ws.Range("H79").Value = 14688.909
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 60000
$ws.Range("N79").Value = -62652
$ws.Range("H115").Value = 2662.3333
$ws.Range("I115").Value = 2498.5
$ws.Range("K115").Value = 7495.5
$ws.Range("M115").Value = -6320.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 535.7
$ws.Range("I97").Value = 419.38095
$ws.Range("J97").Value = 807.1111
$ws.Range("K97").Value = 419.38095
$ws.Range("L97").Value = 807.1111
$ws.Range("M97").Value = 76.61905000000002
$ws.Range("N97").Value = -1799.1111
$ws.Range("H102").Value = 3364.6785
$ws.Range("I102").Value = 2065
$ws.Range("K102").Value = 2065
$ws.Range("M102").Value = -443
$ws.Range("H113").Value = 1765.9286
$ws.Range("I113").Value = 1711.3684
$ws.Range("J113").Value = 1881.1111
$ws.Range("K113").Value = 1711.3684
$ws.Range("L113").Value = 1881.1111
$ws.Range("M113").Value = 458.6315999999999
$ws.Range("N113").Value = -6221.1111
$ws.Range("H122").Value = 2580.0454
$ws.Range("I122").Value = 1302.7894
$ws.Range("J122").Value = 10669.333
$ws.Range("K122").Value = 3908.3682
$ws.Range("L122").Value = 32007.999
$ws.Range("M122").Value = -1458.3682
$ws.Range("N122").Value = -36907.999
$ws.Range("H126").Value = 4425
$ws.Range("I126").Value = 3664.7727
$ws.Range("K126").Value = 10994.3181
$ws.Range("M126").Value = -8524.3181
$ws.Range("H132").Value = 4248.231
$ws.Range("I132").Value = 2118.5
$ws.Range("J132").Value = 7655.8
$ws.Range("K132").Value = 6355.5
$ws.Range("L132").Value = 22967.4
$ws.Range("M132").Value = -3825.5
$ws.Range("N132").Value = -28027.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3881.08
$ws.Range("I16").Value = 3488.3157
$ws.Range("K16").Value = 3488.3157
$ws.Range("M16").Value = -3318.3157
$ws.Range("H22").Value = 1393.68
$ws.Range("I22").Value = 666.2
$ws.Range("J22").Value = 2484.9
$ws.Range("K22").Value = 666.2
$ws.Range("L22").Value = 2484.9
$ws.Range("M22").Value = -371.2
$ws.Range("N22").Value = -3074.9
$ws.Range("H27").Value = 1393.68
$ws.Range("I27").Value = 666.2
$ws.Range("J27").Value = 2484.9
$ws.Range("K27").Value = 666.2
$ws.Range("L27").Value = 2484.9
$ws.Range("M27").Value = -559.2
$ws.Range("N27").Value = -2698.9
$ws.Range("H31").Value = 3269.3333
$ws.Range("I31").Value = 1081.3334
$ws.Range("J31").Value = 9833.333000000001
$ws.Range("K31").Value = 1081.3334
$ws.Range("L31").Value = 9833.333000000001
$ws.Range("M31").Value = -833.3334
$ws.Range("N31").Value = -10329.333
$ws.Range("H46").Value = 1985.0312
$ws.Range("I46").Value = 824.9545000000001
$ws.Range("J46").Value = 4537.2
$ws.Range("K46").Value = 824.9545000000001
$ws.Range("L46").Value = 4537.2
$ws.Range("M46").Value = -636.9545000000001
$ws.Range("N46").Value = -4913.2
$ws.Range("H61").Value = 1931.4193
$ws.Range("I61").Value = 1964.7858
$ws.Range("K61").Value = 1964.7858
$ws.Range("M61").Value = -1762.7858
$ws.Range("H68").Value = 8943.556
$ws.Range("I68").Value = 9641.714
$ws.Range("J68").Value = 6500
$ws.Range("K68").Value = 9641.714
$ws.Range("L68").Value = 6500
$ws.Range("M68").Value = -8892.714
$ws.Range("N68").Value = -7998
$ws.Range("H71").Value = 8943.556
$ws.Range("I71").Value = 9641.714
$ws.Range("J71").Value = 6500
$ws.Range("K71").Value = 48208.57
$ws.Range("L71").Value = 32500
$ws.Range("M71").Value = -44464.57
$ws.Range("N71").Value = -39988
$ws.Range("H80").Value = 45000
$ws.Range("I80").Value = 15000
$ws.Range("J80").Value = 75000
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 75000
$ws.Range("M80").Value = -13877
$ws.Range("N80").Value = -77246
$ws.Range("H83").Value = 45000
$ws.Range("I83").Value = 15000
$ws.Range("J83").Value = 75000
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 225000
$ws.Range("M83").Value = -39384
$ws.Range("N83").Value = -236232
$ws.Range("H93").Value = 939.1818
$ws.Range("I93").Value = 748.25
$ws.Range("K93").Value = 748.25
$ws.Range("M93").Value = 499.75
$ws.Range("H107").Value = 8490
$ws.Range("I107").Value = 8490
$ws.Range("K107").Value = 8490
$ws.Range("M107").Value = -6570
$ws.Range("H113").Value = 1931.4193
$ws.Range("I113").Value = 1964.7858
$ws.Range("K113").Value = 1964.7858
$ws.Range("M113").Value = 205.2141999999999
$ws.Range("H122").Value = 9212
$ws.Range("I122").Value = 5624.5
$ws.Range("K122").Value = 16873.5
$ws.Range("M122").Value = -14423.5
$ws.Range("H132").Value = 1801
$ws.Range("I132").Value = 1095.5883
$ws.Range("K132").Value = 3286.7649
$ws.Range("M132").Value = -756.7648999999997
$ws.Range("H136").Value = 3849.16
$ws.Range("I136").Value = 3188.611
$ws.Range("K136").Value = 9565.832999999999
$ws.Range("M136").Value = -7015.832999999999
$ws.Range("H140").Value = 63087.75
$ws.Range("J140").Value = 63087.75
$ws.Range("L140").Value = 63087.75
$ws.Range("N140").Value = -73447.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 14633.333
$ws.Range("J4").Value = 6950
$ws.Range("L4").Value = 6950
$ws.Range("N4").Value = -7176
$ws.Range("H33").Value = 500
$ws.Range("I33").Value = 500
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 500
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -250
$ws.Range("N33").Value = $null
$ws.Range("H36").Value = 500
$ws.Range("I36").Value = 500
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 500
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -250
$ws.Range("N36").Value = $null
$ws.Range("H51").Value = 39000
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 75000
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 75000
$ws.Range("M51").Value = -2490
$ws.Range("N51").Value = -76020
$ws.Range("H96").Value = 3232.5833
$ws.Range("I96").Value = 3833.3333
$ws.Range("J96").Value = 3032.3333
$ws.Range("K96").Value = 3833.3333
$ws.Range("L96").Value = 3032.3333
$ws.Range("M96").Value = -2460.3333
$ws.Range("N96").Value = -5778.3333
$ws.Range("H99").Value = 48998.332
$ws.Range("J99").Value = 48497.5
$ws.Range("L99").Value = 48497.5
$ws.Range("N99").Value = -54487.5
$ws.Range("H100").Value = 1255.75
$ws.Range("I100").Value = 850.3333
$ws.Range("J100").Value = 1499
$ws.Range("K100").Value = 1700.6666
$ws.Range("L100").Value = 2998
$ws.Range("M100").Value = -1159.6666
$ws.Range("N100").Value = -4080
$ws.Range("H113").Value = 529.11536
$ws.Range("I113").Value = 511.2
$ws.Range("K113").Value = 1533.6
$ws.Range("M113").Value = 636.4000000000001
$ws.Range("H126").Value = 2925.0908
$ws.Range("I126").Value = 2817.7
$ws.Range("J126").Value = 3999
$ws.Range("K126").Value = 8453.099999999999
$ws.Range("L126").Value = 11997
$ws.Range("M126").Value = -5983.099999999999
$ws.Range("N126").Value = -16937
$ws.Range("H132").Value = 2616.9302
$ws.Range("I132").Value = 1816.5526
$ws.Range("K132").Value = 5449.6578
$ws.Range("M132").Value = -2919.6578
$ws.Range("H138").Value = 122994.5
$ws.Range("J138").Value = 122994.5
$ws.Range("L138").Value = 122994.5
$ws.Range("N138").Value = -133274.5
